$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8712884
$ws.Range("I32").Value = 775
$ws.Range("J32").Value = 17424992
$ws.Range("K32").Value = 775
$ws.Range("L32").Value = 17424992
$ws.Range("M32").Value = -449
$ws.Range("N32").Value = -17425644

$ws.Range("H64").Value = 3431.3777
$ws.Range("I64").Value = 2926.5557
$ws.Range("J64").Value = 4188.6113
$ws.Range("K64").Value = 2926.5557
$ws.Range("L64").Value = 4188.6113
$ws.Range("M64").Value = -2678.5557
$ws.Range("N64").Value = -4684.6113

$ws.Range("H67").Value = 3431.3777
$ws.Range("I67").Value = 2926.5557
$ws.Range("J67").Value = 4188.6113
$ws.Range("K67").Value = 2926.5557
$ws.Range("L67").Value = 4188.6113
$ws.Range("M67").Value = -2068.5557
$ws.Range("N67").Value = -5904.6113

$ws.Range("H75").Value = 26462.8
$ws.Range("J75").Value = 26462.8
$ws.Range("L75").Value = 26462.8
$ws.Range("N75").Value = -28334.8

$ws.Range("H78").Value = 26462.8
$ws.Range("J78").Value = 26462.8
$ws.Range("L78").Value = 79388.39999999999
$ws.Range("N78").Value = -88748.39999999999

$ws.Range("H98").Value = 671.3823
$ws.Range("I98").Value = 458.42856
$ws.Range("J98").Value = 1665.1666
$ws.Range("K98").Value = 458.42856
$ws.Range("L98").Value = 1665.1666
$ws.Range("M98").Value = 1039.57144
$ws.Range("N98").Value = -4661.1666

$ws.Range("H122").Value = 671.3823
$ws.Range("I122").Value = 458.42856
$ws.Range("J122").Value = 1665.1666
$ws.Range("K122").Value = 1375.28568
$ws.Range("L122").Value = 4995.4998
$ws.Range("M122").Value = 1074.71432
$ws.Range("N122").Value = -9895.4998

$ws.Range("H136").Value = 38137.6
$ws.Range("J136").Value = 38137.6
$ws.Range("L136").Value = 38137.6
$ws.Range("N136").Value = -48337.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H122").Value = 2444.6667
$ws.Range("I122").Value = 2615.5454
$ws.Range("J122").Value = 2256.7
$ws.Range("K122").Value = 7846.6362
$ws.Range("L122").Value = 6770.099999999999
$ws.Range("M122").Value = -5396.6362
$ws.Range("N122").Value = -11670.1

$ws.Range("H132").Value = 3043.6038
$ws.Range("I132").Value = 2774.838
$ws.Range("J132").Value = 3665.125
$ws.Range("K132").Value = 8324.514000000001
$ws.Range("L132").Value = 10995.375
$ws.Range("M132").Value = -5794.514000000001
$ws.Range("N132").Value = -16055.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7769.5884
$ws.Range("I86").Value = 8900.857
$ws.Range("J86").Value = 6977.7
$ws.Range("K86").Value = 8900.857
$ws.Range("L86").Value = 6977.7
$ws.Range("M86").Value = -7777.857
$ws.Range("N86").Value = -9223.700000000001

$ws.Range("H89").Value = 7769.5884
$ws.Range("I89").Value = 8900.857
$ws.Range("J89").Value = 6977.7
$ws.Range("K89").Value = 44504.285
$ws.Range("L89").Value = 34888.5
$ws.Range("M89").Value = -38888.285
$ws.Range("N89").Value = -46120.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 100668240
$ws.Range("I2").Value = 850
$ws.Range("J2").Value = 302003000
$ws.Range("K2").Value = 850
$ws.Range("L2").Value = 302003000
$ws.Range("M2").Value = -737
$ws.Range("N2").Value = -302003226

$ws.Range("H5").Value = 414.64285
$ws.Range("I5").Value = 79.625
$ws.Range("J5").Value = 861.3333
$ws.Range("K5").Value = 79.625
$ws.Range("L5").Value = 861.3333
$ws.Range("M5").Value = 32.375
$ws.Range("N5").Value = -1085.3333

$ws.Range("H6").Value = 22333594
$ws.Range("I6").Value = 173.66667
$ws.Range("J6").Value = 67000436
$ws.Range("K6").Value = 173.66667
$ws.Range("L6").Value = 67000436
$ws.Range("M6").Value = -60.66667000000001
$ws.Range("N6").Value = -67000662

$ws.Range("H10").Value = 1002.75
$ws.Range("I10").Value = 337
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 337
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -198
$ws.Range("N10").Value = -3278

$ws.Range("H17").Value = 2520
$ws.Range("J17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("N17").Value = -3348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 930
$ws.Range("I80").Value = 825
$ws.Range("K80").Value = 2475
$ws.Range("M80").Value = -1539

$ws.Range("H83").Value = 930
$ws.Range("I83").Value = 825
$ws.Range("K83").Value = 7425
$ws.Range("M83").Value = -2745

$ws.Range("H136").Value = 4217.154
$ws.Range("I136").Value = 4115
$ws.Range("J136").Value = 4304.7144
$ws.Range("K136").Value = 12345
$ws.Range("L136").Value = 12914.1432
$ws.Range("M136").Value = -7245
$ws.Range("N136").Value = -23114.1432

$ws.Range("H138").Value = 2305.6875
$ws.Range("I138").Value = 2017.7778
$ws.Range("K138").Value = 6053.3334
$ws.Range("M138").Value = -913.3334000000004

$ws.Range("H139").Value = 1817.2727
$ws.Range("I139").Value = 838
$ws.Range("J139").Value = 2633.3333
$ws.Range("K139").Value = 2514
$ws.Range("L139").Value = 7899.999899999999
$ws.Range("M139").Value = 2626
$ws.Range("N139").Value = -18179.9999

$ws.Range("H141").Value = 4541.0347
$ws.Range("I141").Value = 3612.8572
$ws.Range("J141").Value = 4836.364
$ws.Range("K141").Value = 10838.5716
$ws.Range("L141").Value = 14509.092
$ws.Range("M141").Value = -5658.571599999999
$ws.Range("N141").Value = -24869.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3049.0833
$ws.Range("I102").Value = 3301.3333
$ws.Range("J102").Value = 1283.3334
$ws.Range("K102").Value = 3301.3333
$ws.Range("L102").Value = 1283.3334
$ws.Range("M102").Value = -1679.3333
$ws.Range("N102").Value = -4527.3334

$ws.Range("H126").Value = 2348.1904
$ws.Range("I126").Value = 1843.6666
$ws.Range("J126").Value = 3020.889
$ws.Range("K126").Value = 5530.9998
$ws.Range("L126").Value = 9062.667000000001
$ws.Range("M126").Value = -3060.9998
$ws.Range("N126").Value = -14002.667

$ws.Range("H132").Value = 4877.033
$ws.Range("I132").Value = 7255.5454
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 21766.6362
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -19236.6362
$ws.Range("N132").Value = -15560

$ws.Range("H141").Value = 36988.6
$ws.Range("J141").Value = 42597
$ws.Range("L141").Value = 42597
$ws.Range("N141").Value = -52957

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1801.2413
$ws.Range("I16").Value = 1777.2916
$ws.Range("J16").Value = 1916.2
$ws.Range("K16").Value = 1777.2916
$ws.Range("L16").Value = 1916.2
$ws.Range("M16").Value = -1607.2916
$ws.Range("N16").Value = -2256.2

$ws.Range("H22").Value = 1407.091
$ws.Range("I22").Value = 1640
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 1640
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -1345
$ws.Range("N22").Value = -1589.5

$ws.Range("H27").Value = 1407.091
$ws.Range("I27").Value = 1640
$ws.Range("J27").Value = 999.5
$ws.Range("K27").Value = 1640
$ws.Range("L27").Value = 999.5
$ws.Range("M27").Value = -1533
$ws.Range("N27").Value = -1213.5

$ws.Range("H41").Value = 12566.667
$ws.Range("J41").Value = 12566.667
$ws.Range("L41").Value = 12566.667
$ws.Range("N41").Value = -13442.667

$ws.Range("H132").Value = 7656.8945
$ws.Range("I132").Value = 2559.3914
$ws.Range("K132").Value = 7678.174199999999
$ws.Range("M132").Value = -5148.174199999999

$ws.Range("H139").Value = 39712.5
$ws.Range("J139").Value = 39712.5
$ws.Range("L139").Value = 39712.5
$ws.Range("N139").Value = -49992.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1808.2616
$ws.Range("I132").Value = 841.86664
$ws.Range("J132").Value = 2636.6
$ws.Range("K132").Value = 2525.59992
$ws.Range("L132").Value = 7909.799999999999
$ws.Range("M132").Value = 4.400080000000344
$ws.Range("N132").Value = -12969.8

$ws.Range("H139").Value = 39995
$ws.Range("J139").Value = 39995
$ws.Range("L139").Value = 39995
$ws.Range("N139").Value = -50275

$ws.Range("H140").Value = 26464.125
$ws.Range("J140").Value = 26464.125
$ws.Range("L140").Value = 26464.125
$ws.Range("N140").Value = -36824.125

$ws.Range("H141").Value = 31810
$ws.Range("J141").Value = 31810
$ws.Range("L141").Value = 31810
$ws.Range("N141").Value = -42170
